# The "Förändrad" (Changed) date column C, for rows 2 through 98,
# was bumped by one day (serial date 45179 -> 45180, i.e. 2023-09-10 -> 2023-09-11).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C98").Value = 45180
